# Diagnoses sheet: keep only the first patient's record, update their
# temperature/symptoms/tally columns, and drop the other three rows
# (Chad's duplicate, Jane Doe, and Cha O) entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Chad Oliver) updates
$ws.Range("E2").Value = 113
$ws.Range("G2").Value = "Difficulty Breathing,Chest Pain,Loss of Movement,Fever,Tiredness,Pains,Sore Throat"
$ws.Range("J2").Value = 3
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 2
$ws.Range("Q2").Value = "Very High Risk"

# Remove rows 3-5 completely (Chad duplicate, Jane Doe, Cha O), shifting
# everything below up and shrinking the used range to A1:Q2.
$ws.Rows("3:5").Delete()
